$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AGALSIDASE BETA"
$ws.Range("B2").Value = "AGALSIDASE BETA  [AGALSIDASE BETA BIOSIMILAR 1]"
$ws.Range("B3").Value = "ALENDRONATE SODIUM HYDRATE"
$ws.Range("A3").Value = "ALENDRONATE SODIUM"
$ws.Range("B4").Value = "AMLODIPINE BESILATE"
$ws.Range("A4").Value = "AMLODIPINE BESYLATE"
$ws.Range("B5").Value = "AMLODIPINE BESILATE; ATORVASTATIN CALCIUM HYDRATE"
$ws.Range("A5").Value = "AMLODIPINE BESYLATE; ATORVASTATIN CALCIUM"
$ws.Range("A6").Value = "AMLODIPINE BESYLATE; VALSARTAN"
$ws.Range("B6").Value = "AMLODIPINE, VALSARTAN"
$ws.Range("A7").Value = "ANAGRELIDE HYDROCHLORIDE"
$ws.Range("B7").Value = "ANAGRELIDE HYDROCHLORIDE HYDRATE"
$ws.Range("A8").Value = "APOMORPHINE HYDROCHLORIDE"
$ws.Range("B8").Value = "APOMORPHINE HYDROCHLORIDE HYDRATE"
$ws.Range("A9").Value = "ARGATROBAN"
$ws.Range("B9").Value = "ARGATROBAN HYDRATE"
$ws.Range("B10").Value = "BELUMOSUDIL MESILATE"
$ws.Range("A10").Value = "BELUMOSUDIL MESYLATE"
$ws.Range("A11").Value = "BETAINE"
$ws.Range("B11").Value = "BETAINE ANHYDROUS"
$ws.Range("A12").Value = "BOSENTAN"
$ws.Range("B12").Value = "BOSENTAN HYDRATE"
$ws.Range("A13").Value = "BOSENTAN"
$ws.Range("B13").Value = "BOSENTAN MONOHYDRATE"
$ws.Range("A14").Value = "BOSUTINIB"
$ws.Range("B14").Value = "BOSUTINIB HYDRATE"
$ws.Range("A15").Value = "BOSUTINIB"
$ws.Range("B15").Value = "BOSUTINIB MONOHYDRATE"
$ws.Range("A16").Value = "BUDESONIDE; FORMOTEROL FUMARATE"
$ws.Range("B16").Value = "BUDESONIDE; FORMOTEROL FUMARATE DIHYDRATE"
$ws.Range("A17").Value = "BUDESONIDE; FORMOTEROL FUMARATE"
$ws.Range("B17").Value = "BUDESONIDE; FORMOTEROL FUMARATE HYDRATE"
$ws.Range("A18").Value = "CANAGLIFLOZIN"
$ws.Range("B18").Value = "CANAGLIFLOZIN HYDRATE"
$ws.Range("A19").Value = "CAPMATINIB HYDROCHLORIDE"
$ws.Range("B19").Value = "CAPMATINIB HYDROCHLORIDE HYDRATE"
$ws.Range("A20").Value = "CASIRIVIMAB, IMDEVIMAB"
$ws.Range("B20").Value = "CASIRIVIMAB; IMDEVIMAB"
$ws.Range("A21").Value = "CEFEPIME DIHYDROCHLORIDE"
$ws.Range("B21").Value = "CEFEPIME HYDROCHLORIDE"
$ws.Range("A22").Value = "CEFIDEROCOL SULFATE TOSYLATE"
$ws.Range("B22").Value = "CEFIDEROCOL TOSILATE SULFATE HYDRATE"

$ws.Range("B22").Select()
